$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.441.64'
$ws.Range('E2').Value = '  +1.78%  '
$ws.Range('D3').Value = '1.854.44'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.60'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6948'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3064'
$ws.Range('E8').Value = '  +0.49%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07669'
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('E10').Value = '  +0.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07777'
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.143'
$ws.Range('E12').Value = '  +1.32%  '
$ws.Range('D13').Value = '1.853.23'
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '91.03'
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6907'
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.277'
$ws.Range('E16').Value = '  -2.25%  '
$ws.Range('D17').Value = '29.429.05'
$ws.Range('E17').Value = '  +1.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008316'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').Value = '2.099.07'
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '237.92'
$ws.Range('E20').Value = '  -2.10%  '
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.603'
$ws.Range('E23').Value = '  +1.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1490'
$ws.Range('E25').Value = '  +1.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.86'
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('E27').Value = '  +0.91%  '
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.531'
$ws.Range('E29').Value = '  -0.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.238'
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.147'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.213'
$ws.Range('E32').Value = '  +2.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05092'
$ws.Range('E33').Value = '  -0.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7712'
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.879'
$ws.Range('E35').Value = '  +2.35%  '
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.679'
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('D38').Value = '1.329.96'
$ws.Range('E38').Value = '  +8.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01869'
$ws.Range('E39').Value = '  +1.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.720'
$ws.Range('E40').Value = '  +0.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9501'
$ws.Range('E41').Value = '  +1.83%  '
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.801'
$ws.Range('E43').Value = '  +2.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.810'
$ws.Range('E45').Value = '  +2.56%  '
$ws.Range('E46').Value = '  +2.18%  '
$ws.Range('D47').Value = '1.997.88'
$ws.Range('E47').Value = '  +0.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5218'
$ws.Range('E48').Value = '  +0.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.781'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '63.07'
$ws.Range('E50').Value = '  -1.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.954'
$ws.Range('E51').Value = '  +0.68%  '
